$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 692.3570999999999
$ws.Range("I12").Value = 474.5
$ws.Range("K12").Value = 474.5
$ws.Range("M12").Value = -304.5
$ws.Range("H18").Value = 2457.1765
$ws.Range("I18").Value = 2392.2666
$ws.Range("K18").Value = 2392.2666
$ws.Range("M18").Value = -2108.2666
$ws.Range("H62").Value = 47772.727
$ws.Range("J62").Value = 47550
$ws.Range("L62").Value = 47550
$ws.Range("N62").Value = -48798
$ws.Range("H65").Value = 47772.727
$ws.Range("J65").Value = 47550
$ws.Range("L65").Value = 237750
$ws.Range("N65").Value = -243990
$ws.Range("H100").Value = 15522190
$ws.Range("I100").Value = 14710958
$ws.Range("J100").Value = 18280378
$ws.Range("K100").Value = 14710958
$ws.Range("L100").Value = 18280378
$ws.Range("M100").Value = -14710417
$ws.Range("N100").Value = -18281460
$ws.Range("H113").Value = 18609
$ws.Range("I113").Value = 26758.285
$ws.Range("J113").Value = 7200
$ws.Range("K113").Value = 26758.285
$ws.Range("L113").Value = 7200
$ws.Range("M113").Value = -23504.285
$ws.Range("N113").Value = -13708
$ws.Range("H137").Value = 3870353.2
$ws.Range("I137").Value = 645244.0600000001
$ws.Range("K137").Value = 1935732.18
$ws.Range("M137").Value = -1933182.18
$ws.Range("H138").Value = 7917.5
$ws.Range("I138").Value = 6033.5884
$ws.Range("K138").Value = 18100.7652
$ws.Range("M138").Value = -12960.7652

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2457.3704
$ws.Range("I32").Value = 2366.8235
$ws.Range("K32").Value = 2366.8235
$ws.Range("M32").Value = -2079.8235
$ws.Range("H44").Value = 84987.75
$ws.Range("J44").Value = 84987.75
$ws.Range("L44").Value = 84987.75
$ws.Range("N44").Value = -85963.75
$ws.Range("H45").Value = 167275.61
$ws.Range("I45").Value = 269596.25
$ws.Range("J45").Value = 3562.6
$ws.Range("K45").Value = 269596.25
$ws.Range("L45").Value = 3562.6
$ws.Range("M45").Value = -269219.25
$ws.Range("N45").Value = -4316.6
$ws.Range("H132").Value = 4618.25
$ws.Range("I132").Value = 4361.077
$ws.Range("K132").Value = 13083.231
$ws.Range("M132").Value = -10553.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 3750
$ws.Range("I24").Value = 3750
$ws.Range("K24").Value = 3750
$ws.Range("M24").Value = -3515
$ws.Range("H134").Value = 9012.625
$ws.Range("I134").Value = 9822.65
$ws.Range("K134").Value = 29467.95
$ws.Range("M134").Value = -26932.95
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 18999.5
$ws.Range("I23").Value = 18999.5
$ws.Range("K23").Value = 18999.5
$ws.Range("M23").Value = -18759.5
$ws.Range("H27").Value = 18999.5
$ws.Range("I27").Value = 18999.5
$ws.Range("K27").Value = 18999.5
$ws.Range("M27").Value = -18807.5
$ws.Range("H31").Value = 2751.0444
$ws.Range("I31").Value = 2071.4119
$ws.Range("K31").Value = 2071.4119
$ws.Range("M31").Value = -1776.4119
$ws.Range("H34").Value = 2751.0444
$ws.Range("I34").Value = 2071.4119
$ws.Range("K34").Value = 2071.4119
$ws.Range("M34").Value = -1869.4119
$ws.Range("H58").Value = 6475.394
$ws.Range("I58").Value = 8851.444
$ws.Range("J58").Value = 3624.1333
$ws.Range("K58").Value = 8851.444
$ws.Range("L58").Value = 3624.1333
$ws.Range("M58").Value = -8648.444
$ws.Range("N58").Value = -4030.1333
$ws.Range("H68").Value = 33089.332
$ws.Range("J68").Value = 44500
$ws.Range("L68").Value = 44500
$ws.Range("N68").Value = -45998
$ws.Range("H71").Value = 33089.332
$ws.Range("J71").Value = 44500
$ws.Range("L71").Value = 133500
$ws.Range("N71").Value = -140988
$ws.Range("H99").Value = 457138.62
$ws.Range("I99").Value = 502352.5
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 502352.5
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -500854.5
$ws.Range("N99").Value = -7996
$ws.Range("H104").Value = 36142
$ws.Range("J104").Value = 36142
$ws.Range("L104").Value = 36142
$ws.Range("N104").Value = -41384
$ws.Range("H126").Value = 457138.62
$ws.Range("I126").Value = 502352.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 1507057.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -1504587.5
$ws.Range("N126").Value = -19940
$ws.Range("H136").Value = 6475.394
$ws.Range("I136").Value = 8851.444
$ws.Range("J136").Value = 3624.1333
$ws.Range("K136").Value = 26554.332
$ws.Range("L136").Value = 10872.3999
$ws.Range("M136").Value = -24004.332
$ws.Range("N136").Value = -15972.3999
$ws.Range("H141").Value = 57377.6
$ws.Range("I141").Value = 50296
$ws.Range("J141").Value = 68000
$ws.Range("K141").Value = 50296
$ws.Range("L141").Value = 68000
$ws.Range("M141").Value = -45116
$ws.Range("N141").Value = -78360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1087.091
$ws.Range("J113").Value = 1373.8096
$ws.Range("L113").Value = 4121.4288
$ws.Range("N113").Value = -8461.4288
$ws.Range("H133").Value = 11882.667
$ws.Range("I133").Value = 5650
$ws.Range("K133").Value = 16950
$ws.Range("M133").Value = -11890

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2500
$ws.Range("I29").Value = 2500
$ws.Range("K29").Value = 2500
$ws.Range("M29").Value = -2210
$ws.Range("H132").Value = 5935.0835
$ws.Range("I132").Value = 3988.6
$ws.Range("K132").Value = 11965.8
$ws.Range("M132").Value = -9435.799999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 42439.184
$ws.Range("I7").Value = 45708.1
$ws.Range("J7").Value = 9750
$ws.Range("K7").Value = 45708.1
$ws.Range("L7").Value = 9750
$ws.Range("M7").Value = -45596.1
$ws.Range("N7").Value = -9974
$ws.Range("H38").Value = 9030
$ws.Range("I38").Value = 9030
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 9030
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -8620
$ws.Range("N38").ClearContents()
$ws.Range("H44").Value = 16000
$ws.Range("J44").Value = 16000
$ws.Range("L44").Value = 16000
$ws.Range("N44").Value = -16912
$ws.Range("H56").Value = 10549.75
$ws.Range("I56").Value = 10549.75
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 10549.75
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -9858.75
$ws.Range("N56").ClearContents()
$ws.Range("H126").Value = 42439.184
$ws.Range("I126").Value = 45708.1
$ws.Range("J126").Value = 9750
$ws.Range("K126").Value = 137124.3
$ws.Range("L126").Value = 29250
$ws.Range("M126").Value = -134654.3
$ws.Range("N126").Value = -34190
$ws.Range("H127").Value = 200297000
$ws.Range("J127").Value = 371248.75
$ws.Range("L127").Value = 371248.75
$ws.Range("N127").Value = -381168.75
$ws.Range("H130").Value = 49249
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 80450
$ws.Range("J131").Value = 80450
$ws.Range("L131").Value = 80450
$ws.Range("N131").Value = -90530
$ws.Range("H136").Value = 4037.077
$ws.Range("I136").Value = 1996.9
$ws.Range("J136").Value = 5312.1875
$ws.Range("K136").Value = 5990.700000000001
$ws.Range("L136").Value = 15936.5625
$ws.Range("M136").Value = -3440.700000000001
$ws.Range("N136").Value = -21036.5625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3402714.8
$ws.Range("I3").Value = 3402714.8
$ws.Range("K3").Value = 3402714.8
$ws.Range("M3").Value = -3402600.8
$ws.Range("H11").Value = 2902577.8
$ws.Range("I11").Value = 3264150
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 3264150
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -3264008
$ws.Range("N11").Value = -10284
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()
$ws.Range("H41").Value = 11869.182
$ws.Range("J41").Value = 11256.1
$ws.Range("L41").Value = 11256.1
$ws.Range("N41").Value = -12036.1
$ws.Range("H126").Value = 36528.152
$ws.Range("J126").Value = 7448.125
$ws.Range("L126").Value = 22344.375
$ws.Range("N126").Value = -27284.375
